$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1 = 'big and tall running pants for men'
    2 = 'girls compression knee sleeve'
    3 = 'football pants with pads'
    4 = 'compression padded shorts'
    5 = 'spandex for men pants'
    6 = 'elastic knee pad'
    7 = 'knee pads gel construction'
    8 = 'capri shorts for men'
    9 = 'football padded pants'
    10 = 'girls volleyball kneepads'
    11 = 'black capris men'
    12 = 'womens basketball pants'
    13 = 'baseball material'
    14 = 'youth sports leggings'
    15 = 'protect knee pads'
    16 = 'knees pad'
    17 = 'black youth knee pads'
    18 = 'exercise pads for knees'
    19 = 'gym shorts above knee for men'
    20 = 'knee sleeve basketball youth'
    21 = 'knee pads construction'
    22 = 'work knee pad'
    23 = 'bjj knee'
    24 = 'knee sleeve bjj'
    25 = 'knee pads under'
    26 = 'mens long cycling pants'
    27 = 'mens pad'
    28 = 'boys long baseball pants'
    29 = 'mens basketball gear'
    30 = 'girl sliding shorts'
    31 = 'calf silicone pads'
    32 = 'compression knee sleeve men pair'
    33 = 'girls black softball pants'
    34 = 'hip protector pads'
    35 = 'flexible work pants for men'
    36 = 'knee protection pads'
    37 = 'lightweight work pants for men'
    38 = 'youth baseball pants long'
    39 = 'knees pads work'
    40 = 'tights boys'
    41 = 'mens shorts long below knee'
    42 = 'knee sleeve youth'
    43 = 'snowboarding pants men'
    44 = 'baseball shorts for men'
    45 = 'baseball mens pants'
    46 = 'knee compression sleeve - reduce strain & swelling'
    47 = 'pads men'
    48 = 'basketball sleeve youth leg'
    49 = 'thigh pads football'
    50 = 'compression volleyball'
    51 = 'leggings for mens'
    52 = 'mens yoga pants'
    53 = 'padded football pants'
    54 = 'spandex capris'
    55 = 'water knee hockey'
    56 = 'compression pants sleeves'
    57 = 'knee sleeve padded'
    58 = 'knees pads for construction'
    59 = 'tight capri'
    60 = 'mens baseball compression shorts'
    61 = 'mens running knee compression'
    62 = 'black football leggings'
    63 = 'knee sleeves basketball youth'
    64 = 'cycling knee pads'
    65 = 'construction knee pad'
    66 = 'compression calf leggings'
    67 = 'baseball youth compression sleeve'
    68 = '6 pairs of leggings'
    69 = 'basketball shorts for men pack of 5'
    70 = 'compression pants youth boys'
    71 = 'impact shorts men'
    72 = 'large knee pad'
    73 = 'baseball compression sleeve'
    74 = 'boys sports leggings'
    75 = 'volleyball spandex pack'
    76 = 'baseball pants youth large'
    77 = 'boys paintball pants'
    78 = 'yoga pants mens'
    79 = 'calf tear compression sleeve'
    80 = 'compression shorts men long length'
    81 = 'yoga hand pads'
    82 = 'knee sleeves with padding'
    83 = 'athletic capri leggings'
    84 = 'pants compression men'
    85 = 'basketball padding'
    86 = 'knee pads for men floor work'
    87 = 'youth knee sleeve wrestling'
    88 = 'professional construction knee pads'
    89 = 'youth basketball'
    90 = 'basketball compression knee sleeve'
    91 = 'black softball pants youth girls'
    92 = 'hex gear wash'
    93 = 'knee construction pads'
    94 = 'youth girls softball pants'
    95 = 'compression tight pants'
    96 = 'male workout leggings'
    97 = 'boys xl baseball pants'
    98 = 'thick leggings for men'
    99 = 'knee pads for'
    100 = 'adult football girdle'
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
